$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) for data rows keeps its text formatting so
# numeric-looking values (e.g. "230.64") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @{
    2 = @{ D="41.709.40"; E="  +5.36%  " }
    3 = @{ D="2.222.27"; E="  +3.14%  " }
    4 = @{ E="  +0.06%  " }
    5 = @{ D="230.64"; E="  +1.69%  " }
    6 = @{ E="  +0.76%  " }
    7 = @{ D="60.95"; E="  -2.66%  " }
    8 = @{ E="  +0.12%  " }
    9 = @{ E="  +2.96%  " }
    10 = @{ D="58.88"; E="  +0.84%  " }
    11 = @{ D="0.0889"; E="  +5.78%  " }
    12 = @{ E="  +0.49%  " }
    13 = @{ D="2.554.24"; E="  +3.20%  " }
    14 = @{ D="15.61"; E="  -1.26%  " }
    15 = @{ D="21.68"; E="  +0.04%  " }
    16 = @{ D="0.796"; E="  -0.71%  " }
    17 = @{ D="5.55"; E="  +1.64%  " }
    18 = @{ D="2.229.28"; E="  +3.11%  " }
    19 = @{ D="41.589.05"; E="  +5.15%  " }
    20 = @{ D="72.75"; E="  +1.66%  " }
    21 = @{ D="0.0₃0894"; E="  +5.46%  " }
    22 = @{ D="6.03"; E="  -0.70%  " }
    23 = @{ D="249.64"; E="  +9.92%  " }
    24 = @{ E="  +0.03%  " }
    25 = @{ D="2.39"; E="  +1.73%  " }
    26 = @{ D="2.28"; E="  +0.30%  " }
    27 = @{ E="  +1.95%  " }
    28 = @{ D="167.46"; E="  -1.88%  " }
    29 = @{ E="  +1.51%  " }
    30 = @{ D="19.93"; E="  +1.80%  " }
    31 = @{ D="1.41"; E="  -0.91%  " }
    32 = @{ D="2.62"; E="  -2.31%  " }
    33 = @{ E="  +0.71%  " }
    34 = @{ D="4.91"; E="  +4.49%  " }
    35 = @{ D="4.60"; E="  +0.83%  " }
    36 = @{ D="0.0622"; E="  +1.25%  " }
    37 = @{ D="6.55"; E="  -5.82%  " }
    38 = @{ D="3.66"; E="  -3.36%  " }
    39 = @{ D="2.34"; E="  -1.81%  " }
    40 = @{ E="  +29.59%  " }
    41 = @{ D="0.999"; E="  -0.17%  " }
    42 = @{ E="  +1.19%  " }
    44 = @{ D="8.57"; E="  +9.38%  " }
    45 = @{ D="0.0979"; E="  +6.43%  " }
    46 = @{ B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="98.81"; E="  -4.00%  " }
    47 = @{ B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="1.21"; E="  +1.42%  " }
    48 = @{ D="1.462.73"; E="  -3.25%  " }
    49 = @{ B="InjectiveProtocol"; C="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D="16.44"; E="  -6.62%  " }
    50 = @{ B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="2.80"; E="  -0.07%  " }
    51 = @{ E="  -1.21%  " }
}

$colIndex = @{ "B" = 2; "C" = 3; "D" = 4; "E" = 5 }

foreach ($r in $updates.Keys) {
    $rowData = $updates[$r]
    foreach ($col in $rowData.Keys) {
        $ws.Cells.Item([int]$r, $colIndex[$col]).Value = $rowData[$col]
    }
}

